$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Test")

$nl = [char]10

# Update row 6 (SES-TC-03) with the correct Session test case data,
# in the same order the new shared strings were introduced.
$ws.Range("E6").Value = "SES-TC-03"
$ws.Range("F6").Value = "Session Module"
$ws.Range("G6").Value = "SES-TS-03"
$ws.Range("D6").Value = "Verify restricted page cannot be accessed " + $nl + "using browser back button after logout"
$ws.Range("H6").Value = "Verify restricted page cannot be accessed " + $nl + "via browser back after logout"
$ws.Range("J6").Value = "1. Login " + $nl + "2. Access restricted page " + $nl + "3. Logout " + $nl + "4. Click browser Back button"
$ws.Range("L6").Value = "System prevents access and redirects user to login page"

# D6 and H6 gain word-wrap in the updated layout.
$ws.Range("D6").WrapText = $true
$ws.Range("H6").WrapText = $true

# Update the saved selection to match the author's final cursor position.
[void]$ws.Range("J8").Select()
